$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Row 2
Set-TextValue 'D2' '37.757.75'
Set-TextValue 'E2' '  -0.13%  '

# Row 3
Set-TextValue 'D3' '2.078.00'
Set-TextValue 'E3' '  -0.23%  '

# Row 4
Set-TextValue 'E4' '  -0.05%  '

# Row 5
Set-TextValue 'D5' '232.45'
Set-TextValue 'E5' '  -0.54%  '

# Row 6
Set-TextValue 'E6' '  -0.24%  '

# Row 7
Set-TextValue 'E7' '  +0.03%  '

# Row 8
Set-TextValue 'D8' '58.13'
Set-TextValue 'E8' '  -1.64%  '

# Row 9
Set-TextValue 'D9' '0.393'
Set-TextValue 'E9' '  -0.56%  '

# Row 10
Set-TextValue 'D10' '0.0780'
Set-TextValue 'E10' '  -1.22%  '

# Row 11
Set-TextValue 'E11' '  +0.42%  '

# Row 12
Set-TextValue 'B12' 'Chainlink'
Set-TextValue 'C12' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D12' '14.86'
Set-TextValue 'E12' '  +0.48%  '

# Row 13
Set-TextValue 'B13' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C13' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D13' '2.383.52'
Set-TextValue 'E13' '  -0.18%  '

# Row 14
Set-TextValue 'D14' '21.28'
Set-TextValue 'E14' '  -0.07%  '

# Row 15
Set-TextValue 'E15' '  -1.11%  '

# Row 16
Set-TextValue 'E16' '  +0.09%  '

# Row 17
Set-TextValue 'D17' '2.076.97'
Set-TextValue 'E17' '  +0.06%  '

# Row 18
Set-TextValue 'D18' '37.683.35'
Set-TextValue 'E18' '  -0.09%  '

# Row 19
Set-TextValue 'D19' '6.16'
Set-TextValue 'E19' '  +0.02%  '

# Row 20
Set-TextValue 'D20' '70.21'
Set-TextValue 'E20' '  -2.16%  '

# Row 21
Set-TextValue 'E21' '  -2.18%  '

# Row 22
Set-TextValue 'D22' '227.83'
Set-TextValue 'E22' '  -0.16%  '

# Row 23
Set-TextValue 'D23' '0.999'
Set-TextValue 'E23' '  +0.00%  '

# Row 24
Set-TextValue 'D24' '2.39'
Set-TextValue 'E24' '  +0.25%  '

# Row 25
Set-TextValue 'E25' '  -2.26%  '

# Row 26
Set-TextValue 'D26' '9.91'
Set-TextValue 'E26' '  +3.16%  '

# Row 27
Set-TextValue 'D27' '169.72'
Set-TextValue 'E27' '  -1.02%  '

# Row 28
Set-TextValue 'D28' '0.132'
Set-TextValue 'E28' '  -3.25%  '

# Row 29
Set-TextValue 'E29' '  -0.95%  '

# Row 30
Set-TextValue 'D30' '1.38'
Set-TextValue 'E30' '  -3.22%  '

# Row 31
Set-TextValue 'E31' '  +0.26%  '

# Row 32
Set-TextValue 'E32' '  -3.10%  '

# Row 33
Set-TextValue 'D33' '0.0630'
Set-TextValue 'E33' '  -0.80%  '

# Row 34
Set-TextValue 'E34' '  -0.54%  '

# Row 35
Set-TextValue 'E35' '  +1.14%  '

# Row 36
Set-TextValue 'E36' '  +0.64%  '

# Row 37
Set-TextValue 'D37' '3.33'
Set-TextValue 'E37' '  -3.34%  '

# Row 38
Set-TextValue 'E38' '  -0.04%  '

# Row 39
Set-TextValue 'E39' '  -1.76%  '

# Row 40
Set-TextValue 'E40' '  +3.93%  '

# Row 41
Set-TextValue 'B41' 'Cronos'
Set-TextValue 'C41' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D41' '0.0960'
Set-TextValue 'E41' '  -2.30%  '

# Row 42
Set-TextValue 'B42' 'Aave'
Set-TextValue 'C42' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D42' '98.34'
Set-TextValue 'E42' '  -0.84%  '

# Row 43
Set-TextValue 'E43' '  +0.46%  '

# Row 44
Set-TextValue 'D44' '1.489.44'
Set-TextValue 'E44' '  +2.57%  '

# Row 45
Set-TextValue 'E45' '  +2.99%  '

# Row 46
Set-TextValue 'D46' '16.89'
Set-TextValue 'E46' '  -2.30%  '

# Row 47
Set-TextValue 'D47' '4.11'
Set-TextValue 'E47' '  -1.43%  '

# Row 49
Set-TextValue 'D49' '7.29'
Set-TextValue 'E49' '  -1.08%  '

# Row 50
Set-TextValue 'E50' '  -0.77%  '

# Row 51
Set-TextValue 'D51' '2.268.19'
Set-TextValue 'E51' '  -0.29%  '
